# Working on assembly reference classes in source generator and authorization
# -- Update the "Permissions" seed-data sheet: drop the explicit Id column from
#    the generated INSERT statements (CreatedAt now uses getdate() instead of a
#    literal null), add a DELETE/reseed helper formula, rename a handful of
#    permission Codes, switch the page to portrait, move the active selection,
#    and drop the now-unused duplicate "Hyperlink" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Permissions")
if (-not $ws) { $ws = $wb.ActiveSheet }

# 1) Rename the permission Code values (column G, rows 3-10).
$ws.Range("G3").Value  = "ReadRole"
$ws.Range("G4").Value  = "EditRole"
$ws.Range("G5").Value  = "InsertRole"
$ws.Range("G6").Value  = "DeleteRole"
$ws.Range("G7").Value  = "ReadUserExtended"
$ws.Range("G8").Value  = "EditUserExtended"
$ws.Range("G9").Value  = "InsertUserExtended"
$ws.Range("G10").Value = "DeleteUserExtended"

# 2) New helper formula in H2: a DELETE + reseed-identity statement.
$ws.Range("H2").Formula = '=CONCATENATE("delete from ",$A$1,"; dbcc checkident (",$A$1,", reseed, 0);")'

# 3) Rewrite the INSERT-statement formula so it no longer emits the Id column
#    and stamps CreatedAt with getdate() instead of a literal null.
$ws.Range("H3").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,", ",$C$2,", ",$D$2,", ",$E$2,", ",$F$2,", ",$G$2,") values(N''",B3,"'', N''",C3,"'', ",IF(TRIM(D3)<>"","N''"&D3&"''","null"),", ",IF(TRIM(E3)<>"","N''"&E3&"''","null")," , getdate(), N''",G3,"'');")'

$ws.Range("H4:H10").Formula = '=CONCATENATE("insert into ",$A$1,"(",$B$2,", ",$C$2,", ",$D$2,", ",$E$2,", ",$F$2,", ",$G$2,") values(N''",B4,"'', N''",C4,"'', ",IF(TRIM(D4)<>"","N''"&D4&"''","null"),", ",IF(TRIM(E4)<>"","N''"&E4&"''","null")," , getdate(), N''",G4,"'');")'

# 4) Give the Code cells (G3:G10) and the new H2 helper the same cell style
#    Excel applied when it regenerated this range.
$ws.Range("G3:G10").Style = "Normal"
$ws.Range("H2").Style = "Normal"

# 5) Switch the sheet to portrait orientation.
$ws.PageSetup.Orientation = 1

# 6) Move the active selection.
$ws.Activate() | Out-Null
$ws.Range("D10").Select() | Out-Null

# 7) Remove the duplicate "Hyperlink" cell style left over in the style table.
$dupHyperlinkIndex = 0
for ($i = 1; $i -le $wb.Styles.Count; $i++) {
    $style = $wb.Styles.Item($i)
    if ($style.Name -eq "Hyperlink") {
        $dupHyperlinkIndex = $i
    }
}
if ($dupHyperlinkIndex -gt 0) {
    $wb.Styles.Item($dupHyperlinkIndex).Delete() | Out-Null
}
